$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 10
$ws.Cells.Item(86, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(86, 3).Value = "La Araucanía"
$ws.Cells.Item(86, 4).Value = 44763
$ws.Cells.Item(86, 5).Value = 9
$ws.Cells.Item(86, 6).Value = 100112005
$ws.Cells.Item(86, 7).Value = "Puerro"
$ws.Cells.Item(86, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 45
$ws.Cells.Item(86, 11).Value = 17000
$ws.Cells.Item(86, 12).Value = 17000
$ws.Cells.Item(86, 13).Value = 17000
$ws.Cells.Item(86, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(86, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(86, 16).Value = 1417
$ws.Cells.Item(86, 17).Value = 12
$ws.Cells.Item(86, 18).Value = "Hortaliza"
